$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 13 mirroring the structure of existing data rows
$ws.Cells.Item(13, 1).Value = 13
$ws.Cells.Item(13, 2).Value = "delivery"
$ws.Cells.Item(13, 3).Value = "[4, 4, 6, 3, 3, 2]"
$ws.Cells.Item(13, 4).Value = $false
$ws.Cells.Item(13, 5).Value = "InProgress"
$ws.Cells.Item(13, 6).Value = 4
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
